# "new made boss qiongqi and unicorn"
#
# Adds two new monster/boss records to the 表2 table on the "怪物"
# (Monster) worksheet: 穷奇 / qiongqi and 独角兽 / unicorn.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # 怪物 (Monster) sheet
$ws.Activate()

$lo = $ws.ListObjects.Item(1)  # 表2

# Grow the table by two rows (extends the table range + autofilter from
# A3:S43 to A3:S45, same as typing new rows right below the table).
$newRow1 = $lo.ListRows.Add()
$newRow2 = $lo.ListRows.Add()

# Carry the formatting (borders / the yellow "World" column highlight)
# of the previous last data row down into the two freshly added rows so
# the table keeps a consistent look.
$ws.Range("A43:S43").Copy()
$ws.Range("A44:S45").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row 44: 穷奇 / qiongqi
# ---------------------------------------------------------------------
$ws.Cells.Item(44, 1).Value = 43002001    # Id
$ws.Cells.Item(44, 2).Value = "穷奇"       # Name
$ws.Cells.Item(44, 3).Value = "qiongqi"   # Ename
$ws.Cells.Item(44, 4).Value = 0           # Type
$ws.Cells.Item(44, 5).Value = "幻兽传说"   # World
$ws.Cells.Item(44, 6).Value = "atr0"      # Deck
$ws.Cells.Item(44, 7).Value = 11001003    # Job
$ws.Cells.Item(44, 8).Value = 7           # Level
$ws.Cells.Item(44, 10).Value = "common"   # Method
$ws.Cells.Item(44, 11).Value = "qiongqi"  # Emethod
$ws.Cells.Item(44, 12).Value = 1          # EpSlow
$ws.Cells.Item(44, 13).Value = 15         # CardReduce
$ws.Cells.Item(44, 16).Value = "qiongqi"  # Figue
$ws.Cells.Item(44, 17).Value = "oneline"  # BattleMap

# ---------------------------------------------------------------------
# Row 45: 独角兽 / unicorn
# ---------------------------------------------------------------------
$ws.Cells.Item(45, 1).Value = 43002002    # Id
$ws.Cells.Item(45, 2).Value = "独角兽"     # Name
$ws.Cells.Item(45, 3).Value = "unicorn"   # Ename
$ws.Cells.Item(45, 4).Value = 0           # Type
$ws.Cells.Item(45, 5).Value = "幻兽传说"   # World
$ws.Cells.Item(45, 6).Value = "atr0"      # Deck
$ws.Cells.Item(45, 7).Value = 11001003    # Job
$ws.Cells.Item(45, 8).Value = 7           # Level
$ws.Cells.Item(45, 10).Value = "common"   # Method
$ws.Cells.Item(45, 11).Value = "unicorn"  # Emethod
$ws.Cells.Item(45, 12).Value = 1          # EpSlow
$ws.Cells.Item(45, 13).Value = 15         # CardReduce
$ws.Cells.Item(45, 16).Value = "unicorn"  # Figue
$ws.Cells.Item(45, 17).Value = "oneline"  # BattleMap

# Leave the selection on the newly typed Figue column, like a user would
# right after finishing data entry on the new rows.
$ws.Range("P44:P45").Select()

Write-Output "Added qiongqi (row 44) and unicorn (row 45) to the monster table"
